# Auto-generated edit script applying the Tiamat_Profits market-data refresh
# described in the commit "chore: update Sheets via scheduled runner".
# Updates cached price/profit figures (columns H-N) for the affected leve rows
# across the ALC, ARM, BSM, CRP, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# --- ALC row 43 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 67904.39999999999  # H43: was 74553.734
$ws.Cells.Item(43, 9).Value = 0  # I43: was 33506.668
$ws.Cells.Item(43, 10).Value = 67904.39999999999  # J43: was 84815.5
$ws.Cells.Item(43, 11).Value = 0  # K43: was 33506.668
$ws.Cells.Item(43, 12).Value = 67904.39999999999  # L43: was 84815.5
$ws.Cells.Item(43, 13).ClearContents()  # M43: was -33437.668
$ws.Cells.Item(43, 14).Value = -68042.39999999999  # N43: was -84953.5

# --- ALC row 53 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 276  # H53: was 358
$ws.Cells.Item(53, 9).Value = 25.2  # I53: was 26.5
$ws.Cells.Item(53, 10).Value = 694  # J53: was 1021
$ws.Cells.Item(53, 11).Value = 25.2  # K53: was 26.5
$ws.Cells.Item(53, 12).Value = 694  # L53: was 1021
$ws.Cells.Item(53, 13).Value = 611.8  # M53: was 610.5
$ws.Cells.Item(53, 14).Value = -1968  # N53: was -2295

# --- ALC row 55 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 932.7  # H55: was 932.6
$ws.Cells.Item(55, 9).Value = 1307  # I55: was 1709
$ws.Cells.Item(55, 10).Value = 683.1667  # J55: was 599.8570999999999
$ws.Cells.Item(55, 11).Value = 1307  # K55: was 1709
$ws.Cells.Item(55, 12).Value = 683.1667  # L55: was 599.8570999999999
$ws.Cells.Item(55, 13).Value = -1093  # M55: was -1495
$ws.Cells.Item(55, 14).Value = -1111.1667  # N55: was -1027.8571

# --- ALC row 70 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 2355.2144  # H70: was 2964.3333
$ws.Cells.Item(70, 9).Value = 1815  # I70: was 2518
$ws.Cells.Item(70, 10).Value = 2760.375  # J70: was 3283.1428
$ws.Cells.Item(70, 11).Value = 5445  # K70: was 7554
$ws.Cells.Item(70, 12).Value = 8281.125  # L70: was 9849.428400000001
$ws.Cells.Item(70, 13).Value = -5175  # M70: was -7284
$ws.Cells.Item(70, 14).Value = -8821.125  # N70: was -10389.4284

# --- ALC row 73 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 2355.2144  # H73: was 2964.3333
$ws.Cells.Item(73, 9).Value = 1815  # I73: was 2518
$ws.Cells.Item(73, 10).Value = 2760.375  # J73: was 3283.1428
$ws.Cells.Item(73, 11).Value = 5445  # K73: was 7554
$ws.Cells.Item(73, 12).Value = 8281.125  # L73: was 9849.428400000001
$ws.Cells.Item(73, 13).Value = -4509  # M73: was -6618
$ws.Cells.Item(73, 14).Value = -10153.125  # N73: was -11721.4284

# --- ALC row 113 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 1796.069  # H113: was 1782.8695
$ws.Cells.Item(113, 9).Value = 1726.6666  # I113: was 1687.5
$ws.Cells.Item(113, 10).Value = 1909.6364  # J113: was 2000.8572
$ws.Cells.Item(113, 11).Value = 1726.6666  # K113: was 1687.5
$ws.Cells.Item(113, 12).Value = 1909.6364  # L113: was 2000.8572
$ws.Cells.Item(113, 13).Value = 1527.3334  # M113: was 1566.5
$ws.Cells.Item(113, 14).Value = -8417.636399999999  # N113: was -8508.8572

# --- ALC row 116 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 10653.667  # H116: was 10587.4
$ws.Cells.Item(116, 9).Value = 19167.5  # I116: was 16829.285
$ws.Cells.Item(116, 10).Value = 4977.778  # J116: was 5125.75
$ws.Cells.Item(116, 11).Value = 19167.5  # K116: was 16829.285
$ws.Cells.Item(116, 12).Value = 4977.778  # L116: was 5125.75
$ws.Cells.Item(116, 13).Value = -15725.5  # M116: was -13387.285
$ws.Cells.Item(116, 14).Value = -11861.778  # N116: was -12009.75

# --- ARM row 74 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 51095.05  # H74: was 25448.219
$ws.Cells.Item(74, 9).Value = 58114.37  # I74: was 28710.5
$ws.Cells.Item(74, 11).Value = 58114.37  # K74: was 28710.5
$ws.Cells.Item(74, 13).Value = -57240.37  # M74: was -27836.5

# --- ARM row 77 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 51095.05  # H77: was 25448.219
$ws.Cells.Item(77, 9).Value = 58114.37  # I77: was 28710.5
$ws.Cells.Item(77, 11).Value = 290571.85  # K77: was 143552.5
$ws.Cells.Item(77, 13).Value = -286203.85  # M77: was -139184.5

# --- BSM row 134 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 32523.605  # H134: was 40891.73
$ws.Cells.Item(134, 9).Value = 1382.0741  # I134: was 1529.174
$ws.Cells.Item(134, 10).Value = 172660.5  # J134: was 342671.34
$ws.Cells.Item(134, 11).Value = 4146.2223  # K134: was 4587.522
$ws.Cells.Item(134, 12).Value = 517981.5  # L134: was 1028014.02
$ws.Cells.Item(134, 13).Value = -1611.2223  # M134: was -2052.522
$ws.Cells.Item(134, 14).Value = -523051.5  # N134: was -1033084.02

# --- CRP row 16 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 800  # H16: was 749.7143
$ws.Cells.Item(16, 9).Value = 600  # I16: was 671.8
$ws.Cells.Item(16, 10).Value = 1200  # J16: was 793
$ws.Cells.Item(16, 11).Value = 600  # K16: was 671.8
$ws.Cells.Item(16, 12).Value = 1200  # L16: was 793
$ws.Cells.Item(16, 13).Value = -313  # M16: was -384.8
$ws.Cells.Item(16, 14).Value = -1774  # N16: was -1367

# --- CRP row 21 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(21, 8).Value = 4999  # H21: was 35000
$ws.Cells.Item(21, 9).Value = 4999  # I21: was 0
$ws.Cells.Item(21, 10).Value = 0  # J21: was 35000
$ws.Cells.Item(21, 11).Value = 4999  # K21: was 0
$ws.Cells.Item(21, 12).Value = 0  # L21: was 35000
$ws.Cells.Item(21, 13).Value = -4764  # M21: was None
$ws.Cells.Item(21, 14).ClearContents()  # N21: was -35470

# --- CRP row 23 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(23, 8).Value = 70010  # H23: was 58006
$ws.Cells.Item(23, 9).Value = 0  # I23: was 70000
$ws.Cells.Item(23, 10).Value = 70010  # J23: was 55007.5
$ws.Cells.Item(23, 11).Value = 0  # K23: was 70000
$ws.Cells.Item(23, 12).Value = 70010  # L23: was 55007.5
$ws.Cells.Item(23, 13).ClearContents()  # M23: was -69760
$ws.Cells.Item(23, 14).Value = -70490  # N23: was -55487.5

# --- CRP row 27 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(27, 8).Value = 70010  # H27: was 58006
$ws.Cells.Item(27, 9).Value = 0  # I27: was 70000
$ws.Cells.Item(27, 10).Value = 70010  # J27: was 55007.5
$ws.Cells.Item(27, 11).Value = 0  # K27: was 70000
$ws.Cells.Item(27, 12).Value = 70010  # L27: was 55007.5
$ws.Cells.Item(27, 13).ClearContents()  # M27: was -69808
$ws.Cells.Item(27, 14).Value = -70394  # N27: was -55391.5

# --- CRP row 59 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(59, 8).Value = 11024  # H59: was 11624
$ws.Cells.Item(59, 10).Value = 11950.77  # J59: was 12097.143
$ws.Cells.Item(59, 12).Value = 11950.77  # L59: was 12097.143
$ws.Cells.Item(59, 14).Value = -14240.77  # N59: was -14387.143

# --- CRP row 68 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(68, 8).Value = 31836  # H68: was 27454.285
$ws.Cells.Item(68, 9).Value = 0  # I68: was 4000
$ws.Cells.Item(68, 10).Value = 31836  # J68: was 31363.334
$ws.Cells.Item(68, 11).Value = 0  # K68: was 4000
$ws.Cells.Item(68, 12).Value = 31836  # L68: was 31363.334
$ws.Cells.Item(68, 13).ClearContents()  # M68: was -3251
$ws.Cells.Item(68, 14).Value = -33334  # N68: was -32861.334

# --- CRP row 71 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(71, 8).Value = 31836  # H71: was 27454.285
$ws.Cells.Item(71, 9).Value = 0  # I71: was 4000
$ws.Cells.Item(71, 10).Value = 31836  # J71: was 31363.334
$ws.Cells.Item(71, 11).Value = 0  # K71: was 12000
$ws.Cells.Item(71, 12).Value = 95508  # L71: was 94090.00199999999
$ws.Cells.Item(71, 13).ClearContents()  # M71: was -8256
$ws.Cells.Item(71, 14).Value = -102996  # N71: was -101578.002

# --- CRP row 113 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 800  # H113: was 749.7143
$ws.Cells.Item(113, 9).Value = 600  # I113: was 671.8
$ws.Cells.Item(113, 10).Value = 1200  # J113: was 793
$ws.Cells.Item(113, 11).Value = 600  # K113: was 671.8
$ws.Cells.Item(113, 12).Value = 1200  # L113: was 793
$ws.Cells.Item(113, 13).Value = 1570  # M113: was 1498.2
$ws.Cells.Item(113, 14).Value = -5540  # N113: was -5133

# --- GSM row 126 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 1716.4073  # H126: was 1763.6086
$ws.Cells.Item(126, 9).Value = 1430.6428  # I126: was 1480.6923
$ws.Cells.Item(126, 10).Value = 2024.1538  # J126: was 2131.4
$ws.Cells.Item(126, 11).Value = 4291.928400000001  # K126: was 4442.0769
$ws.Cells.Item(126, 12).Value = 6072.4614  # L126: was 6394.200000000001
$ws.Cells.Item(126, 13).Value = -1821.928400000001  # M126: was -1972.0769
$ws.Cells.Item(126, 14).Value = -11012.4614  # N126: was -11334.2

# --- GSM row 132 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 35415.066  # H132: was 62827.293
$ws.Cells.Item(132, 9).Value = 2238.3333  # I132: was 2866.6667
$ws.Cells.Item(132, 10).Value = 66518.25  # J132: was 75676
$ws.Cells.Item(132, 11).Value = 6714.999899999999  # K132: was 8600.000100000001
$ws.Cells.Item(132, 12).Value = 199554.75  # L132: was 227028
$ws.Cells.Item(132, 13).Value = -4184.999899999999  # M132: was -6070.000100000001
$ws.Cells.Item(132, 14).Value = -204614.75  # N132: was -232088

# --- LTW row 7 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4631.643  # H7: was 4357.3125
$ws.Cells.Item(7, 9).Value = 6473.2856  # I7: was 8599.799999999999
$ws.Cells.Item(7, 10).Value = 2790  # J7: was 2428.9092
$ws.Cells.Item(7, 11).Value = 6473.2856  # K7: was 8599.799999999999
$ws.Cells.Item(7, 12).Value = 2790  # L7: was 2428.9092
$ws.Cells.Item(7, 13).Value = -6361.2856  # M7: was -8487.799999999999
$ws.Cells.Item(7, 14).Value = -3014  # N7: was -2652.9092

# --- LTW row 24 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(24, 8).Value = 21652.1  # H24: was 70007
$ws.Cells.Item(24, 9).Value = 1100  # I24: was 0
$ws.Cells.Item(24, 10).Value = 23935.666  # J24: was 70007
$ws.Cells.Item(24, 11).Value = 1100  # K24: was 0
$ws.Cells.Item(24, 12).Value = 23935.666  # L24: was 70007
$ws.Cells.Item(24, 13).Value = -757  # M24: was None
$ws.Cells.Item(24, 14).Value = -24621.666  # N24: was -70693

# --- LTW row 40 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 43805.875  # H40: was 48418.863
$ws.Cells.Item(40, 9).Value = 1707.4  # I40: was 1954.5454
$ws.Cells.Item(40, 10).Value = 113970  # J40: was 94883.17999999999
$ws.Cells.Item(40, 11).Value = 1707.4  # K40: was 1954.5454
$ws.Cells.Item(40, 12).Value = 113970  # L40: was 94883.17999999999
$ws.Cells.Item(40, 13).Value = -1571.4  # M40: was -1818.5454
$ws.Cells.Item(40, 14).Value = -114242  # N40: was -95155.17999999999

# --- LTW row 126 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 4631.643  # H126: was 4357.3125
$ws.Cells.Item(126, 9).Value = 6473.2856  # I126: was 8599.799999999999
$ws.Cells.Item(126, 10).Value = 2790  # J126: was 2428.9092
$ws.Cells.Item(126, 11).Value = 19419.8568  # K126: was 25799.4
$ws.Cells.Item(126, 12).Value = 8370  # L126: was 7286.7276
$ws.Cells.Item(126, 13).Value = -16949.8568  # M126: was -23329.4
$ws.Cells.Item(126, 14).Value = -13310  # N126: was -12226.7276

# --- LTW row 132 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 573899.9399999999  # H132: was 89554.914
$ws.Cells.Item(132, 9).Value = 127137.31  # I132: was 119963.586
$ws.Cells.Item(132, 10).Value = 2003540.4  # J132: was 3397
$ws.Cells.Item(132, 11).Value = 381411.93  # K132: was 359890.758
$ws.Cells.Item(132, 12).Value = 6010621.199999999  # L132: was 10191
$ws.Cells.Item(132, 13).Value = -378881.93  # M132: was -357360.758
$ws.Cells.Item(132, 14).Value = -6015681.199999999  # N132: was -15251

# --- LTW row 136 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2462.2307  # H136: was 2343.7856
$ws.Cells.Item(136, 9).Value = 976.8  # I136: was 948
$ws.Cells.Item(136, 11).Value = 2930.4  # K136: was 2844
$ws.Cells.Item(136, 13).Value = -380.3999999999996  # M136: was -294

# --- WVR row 23 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(23, 8).Value = 28704.2  # H23: was 38005.5
$ws.Cells.Item(23, 9).Value = 1166.3334  # I23: was 9000
$ws.Cells.Item(23, 10).Value = 70011  # J23: was 47674
$ws.Cells.Item(23, 11).Value = 1166.3334  # K23: was 9000
$ws.Cells.Item(23, 12).Value = 70011  # L23: was 47674
$ws.Cells.Item(23, 13).Value = -937.3334  # M23: was -8771
$ws.Cells.Item(23, 14).Value = -70469  # N23: was -48132

# --- WVR row 136 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 2909225.8  # H136: was 2468511.8
$ws.Cells.Item(136, 9).Value = 2646604.2  # I136: was 2464093
$ws.Cells.Item(136, 10).Value = 10000005  # J136: was 2500546.2
$ws.Cells.Item(136, 11).Value = 7939812.600000001  # K136: was 7392279
$ws.Cells.Item(136, 12).Value = 30000015  # L136: was 7501638.600000001
$ws.Cells.Item(136, 13).Value = -7937262.600000001  # M136: was -7389729
$ws.Cells.Item(136, 14).Value = -30005115  # N136: was -7506738.600000001
